$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.514.98"
$ws.Range("E2").Value = "  -0.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.517.48"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.67"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.88"
$ws.Range("E6").Value = "  -2.88%  "

# Row 7
$ws.Range("E7").Value = "  -1.28%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  -2.92%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.22"
$ws.Range("E10").Value = "  -2.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  -0.73%  "

# Row 12
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.20"
$ws.Range("E13").Value = "  -2.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.898.57"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.27"
$ws.Range("E15").Value = "  -3.42%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.514.30"
$ws.Range("E16").Value = "  -1.99%  "

# Row 17
$ws.Range("E17").Value = "  -3.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.496.67"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("E19").Value = "  -2.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0944"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.97"
$ws.Range("E21").Value = "  -3.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.93"
$ws.Range("E22").Value = "  -0.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.94"
$ws.Range("E23").Value = "  -2.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.85"
$ws.Range("E24").Value = "  -3.20%  "

# Row 25
$ws.Range("E25").Value = "  -3.61%  "

# Row 26
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.21"
$ws.Range("E27").Value = "  -4.92%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -4.15%  "

# Row 29
$ws.Range("E29").Value = "  -1.32%  "

# Row 30
$ws.Range("E30").Value = "  -6.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  +1.99%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.60"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.81"
$ws.Range("E33").Value = "  +6.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.66"
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0785"
$ws.Range("E35").Value = "  -2.29%  "

# Row 36
$ws.Range("E36").Value = "  -3.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.43"
$ws.Range("E38").Value = "  -5.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  -3.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("E40").Value = "  -0.89%  "

# Row 41
$ws.Range("E41").Value = "  -0.80%  "

# Row 42
$ws.Range("E42").Value = "  -2.95%  "

# Row 43
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
$ws.Range("E44").Value = "  -1.23%  "

# Row 45
$ws.Range("E45").Value = "  -1.31%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.999.46"
$ws.Range("E46").Value = "  +0.91%  "

# Row 47
$ws.Range("E47").Value = "  +0.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.754.60"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
$ws.Range("E49").Value = "  -2.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.79"
$ws.Range("E50").Value = "  -3.20%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.40"
$ws.Range("E51").Value = "  -2.97%  "
